# previsao_retorno.xlsx - refresh of "Resumo_por_Cliente" report data
# (dados bibi - atualizacao diaria dos calculos de "meses sem comprar"
#  e novo registro de compra para BEMOL S/A)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumo_por_Cliente")

# --- "situacao" (column J) bucket refresh -------------------------------
# The underlying "days since last purchase" metric advanced by one day for
# the whole report, so most INATIVO rows' "X.Y meses sem comprar" label
# ticks up by 0.1.
$ws.Range("J2").Value   = "INATIVO - 54.5 meses sem comprar"
$ws.Range("J8").Value   = "INATIVO - 17.2 meses sem comprar"
$ws.Range("J17").Value  = "INATIVO - 36.6 meses sem comprar"
$ws.Range("J19").Value  = "INATIVO - 14.5 meses sem comprar"
$ws.Range("J45").Value  = "INATIVO - 1.7 meses sem comprar"
$ws.Range("J73").Value  = "INATIVO - 32.6 meses sem comprar"
$ws.Range("J81").Value  = "INATIVO - 25.4 meses sem comprar"
$ws.Range("J86").Value  = "INATIVO - 4.1 meses sem comprar"
$ws.Range("J90").Value  = "INATIVO - 32.6 meses sem comprar"
$ws.Range("J91").Value  = "INATIVO - 12.9 meses sem comprar"
$ws.Range("J92").Value  = "INATIVO - 17.9 meses sem comprar"
$ws.Range("J97").Value  = "INATIVO - 1.4 meses sem comprar"
$ws.Range("J102").Value = "INATIVO - 24.1 meses sem comprar"
$ws.Range("J103").Value = "INATIVO - 9.7 meses sem comprar"
$ws.Range("J104").Value = "INATIVO - 24.5 meses sem comprar"
$ws.Range("J105").Value = "INATIVO - 14.9 meses sem comprar"
$ws.Range("J106").Value = "INATIVO - 5.8 meses sem comprar"
$ws.Range("J110").Value = "INATIVO - 7.7 meses sem comprar"

# --- BEMOL S/A (row 111) - new purchase recorded ------------------------
$ws.Range("E111").Value = 15189
$ws.Range("H111").Value = Get-Date -Year 2025 -Month 6 -Day 3 -Hour 17 -Minute 42 -Second 1
$ws.Range("I111").Value = Get-Date -Year 2025 -Month 6 -Day 4 -Hour 17 -Minute 42 -Second 1
